$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 967, shifting existing rows 967..1058 down to 968..1059
$ws.Rows.Item(967).Insert()

# Populate the newly inserted row 967 with the new record's data
$ws.Range("A967").Value = 3
$ws.Range("B967").Value = "Femacal de La Calera"
$ws.Range("C967").Value = "Coquimbo"
$ws.Range("D967").Value = 45132
$ws.Range("E967").Value = 5
$ws.Range("F967").Value = 100112006
$ws.Range("G967").Value = "Repollo"
$ws.Range("H967").Value = "Crespo record"
$ws.Range("I967").Value = "Primera"
$ws.Range("J967").Value = 3200
$ws.Range("K967").Value = 750
$ws.Range("L967").Value = 800
$ws.Range("M967").Value = 780
$ws.Range("N967").Value = "$/unidad"
$ws.Range("O967").Value = "Provincia de Quillota"
$ws.Range("P967").Value = 780
$ws.Range("Q967").Value = 1
$ws.Range("R967").Value = "Hortaliza"
